$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 17 data (matches existing table columns: A=Question, B=Difficulty,
# C=Pattern, D=Notes, E=Link)
$ws.Range("A17").Value = "2115. Find All Possible Recipes from Given Supplies"
$ws.Range("B17").Value = "Medium"
$ws.Range("C17").Value = "Topological Sort"

# Turn the Link cell into a real hyperlink, matching the style used by the
# other Link cells in column E (Hyperlink cell style). The hyperlink target
# address itself is trimmed (no trailing space), same as the other rows,
# while the displayed cell text keeps the trailing space.
$ws.Hyperlinks.Add($ws.Range("E17"), "https://leetcode.com/problems/find-all-possible-recipes-from-given-supplies/solutions/1646584/java-python-3-toplogical-sort-w-brief-explanation/")
$ws.Range("E17").Value = "https://leetcode.com/problems/find-all-possible-recipes-from-given-supplies/solutions/1646584/java-python-3-toplogical-sort-w-brief-explanation/ "
$ws.Range("E17").Style = "Hyperlink"

$ws.Range("D17").Value = "Brute force is repeated BFS. Optimal is Topological Sort. For each recipe, count of its dependent ingredients as degree, and store (ingredient, recipes that depend on it) as HashMap. Use supplies as the starting points of the topological sort. Use top sort to decrease the in degree of recipes, when in degree reaches 0, add to return list."

# Match the existing "Medium" row styling (orange fill) used by the other
# Medium rows, e.g. B3:B11 / B16.
$ws.Range("B17").Interior.Color = $ws.Range("B16").Interior.Color

# Grow the table (ListObject) to include the new row so the autofilter /
# table range covers A1:E17 instead of A1:E16.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E17"))

# Move the active selection the way Excel would after typing into D17 and
# pressing Enter (cursor advances to D18).
$ws.Range("D18").Select()
